# automate_tests/moisture_data.xlsx
#
# Exercises two helpers added for the test-automation suite against the
# "raw" sheet:
#   - an "append data" routine that (re)writes a 0..30 numeric header row
#     and drops a couple of probe values further down the sheet, and
#   - a "find row" routine that locates the first free row and writes a
#     sentinel value into it.
# The "raw" sheet ends up the active/selected tab, matching the state the
# automation run left the workbook in.

$wb  = $excel.ActiveWorkbook
$raw = $wb.Worksheets.Item("raw")

# --- append data: rewrite header row 1 as a plain numeric 0..30 sequence,
#     clearing the previous "Moisture Content(%)" label + its wrap style ---
$headerRange = $raw.Range("A1:AF1")
$headerRange.Clear()
for ($col = 1; $col -le 31; $col++) {
    $raw.Cells.Item(1, $col).Value = $col - 1
}

# a stray appended value, one row down
$raw.Cells.Item(2, 29).Value = 15

# --- find row: locate the first empty row beneath the data and mark it ---
$targetRow = 11
$raw.Cells.Item($targetRow, 6).Value = " "

# "raw" becomes the active tab/selection after the test run
$raw.Activate()
$raw.Cells.Item($targetRow, 6).Select()
